$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Row 3 corresponds to 63191932-8840-42c6-a8ac-bd4d0e656de1.md
# Status for zh-cn (B3) and de-de (C3) moves from "Handed back: in sync with en-US"
# to "Ready for handoff"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
# Row 3 = 63191932-8840-42c6-a8ac-bd4d0e656de1.md entry: mark as ready for
# handoff again and refresh the handoff timestamp. Row 2 shared the exact
# same handoff timestamp text as row 3, so it is refreshed identically.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-25 06:24:10"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-02-25 06:24:10"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-25 06:24:22"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-02-25 06:24:22"
